# Datum_Conversion.xlsx update: re-order / prune of trigger-station rows for
# the CPRA slide deck generation tools (SlideDeckUpdates_07012020 refresh).
#
# Three stations are no longer used as triggers and their rows are removed
# entirely (cells shift up, shared-string table is recompacted by Excel):
#   - B6  (82770)
#   - B17 (76220)
#   - B21 (76593)
#
# Deleting from the bottom row upward keeps the remaining row numbers stable
# while we work.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(21).EntireRow.Delete()
$ws.Rows.Item(17).EntireRow.Delete()
$ws.Rows.Item(6).EntireRow.Delete()

# Leave the selection where the author's last edit landed.
$ws.Range("C29").Select()
